# Daily attendance processing - 2025-11-25 07:03:19
# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (column G) wherever the value is exactly
# "System, <email>" (a single other recorder paired with System).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value()

    if ($value -ne $null -and $value -is [string] -and $value.StartsWith("System, ")) {
        $rest = $value.Substring(8)
        if (-not $rest.Contains(",")) {
            $cell.Value = "$rest, System"
        }
    }
}
